$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Replace Substrings")
$dst = $wb.Worksheets.Item("Text Case")

$xlPasteFormats = -4122

# ------------------------------------------------------------------
# 1) Shift the existing "Python" table one row down and one column
#    right (A1:C6 -> B2:D7) to make room for the new banner row and
#    the new "Orange" table, using a scratch range on the same sheet
#    so nothing is lost.
# ------------------------------------------------------------------
$dst.Range("A1:C6").Copy($dst.Range("Z1:AB6"))
$dst.Range("A1:C6").Clear()
$dst.Range("Z1:AB5").Copy($dst.Range("B2:D6"))
$dst.Range("Z6:AA6").Copy($dst.Range("B7:C7"))
$dst.Range("Z1:AB6").Clear()

# Column D (the wrapped "Content" column) needs the wrap-text variant
# of the style used in the donor "Replace Substrings" sheet.
$src.Range("D2").Copy()
$dst.Range("D2").PasteSpecial($xlPasteFormats)
$src.Range("D3:D6").Copy()
$dst.Range("D3:D6").PasteSpecial($xlPasteFormats)

# ------------------------------------------------------------------
# 2) Banner row (Python / Orange / Data Polish headers) - identical
#    to the "Replace Substrings" sheet.
# ------------------------------------------------------------------
$src.Range("C1").Copy($dst.Range("C1"))
$src.Range("D1").Copy($dst.Range("D1"))
$src.Range("G1").Copy($dst.Range("G1"))
$src.Range("H1").Copy($dst.Range("H1"))
$src.Range("K1").Copy($dst.Range("K1"))

# ------------------------------------------------------------------
# 3) Right ("Orange") table: header + first/last data rows match the
#    "Replace Substrings" sheet's Orange table verbatim; rows 4-5 get
#    new Text-Case-specific copy (reuse their formatting by copying
#    first, then overwriting just the text).
# ------------------------------------------------------------------
$src.Range("F2:H3").Copy($dst.Range("F2:H3"))
$src.Range("F4:H5").Copy($dst.Range("F4:H5"))
$src.Range("F6:H7").Copy($dst.Range("F6:H7"))

$dst.Range("F4").Value = "Inspect Text Case"
$dst.Range("H4").Value = "Use 'Data Table' to view text case."
$dst.Range("F5").Value = "Change Text Case"
$dst.Range("H5").Value = "Connect 'Preprocess Text', select 'Change Case'."

# ------------------------------------------------------------------
# 4) Column widths / row heights.
# ------------------------------------------------------------------
$dst.Columns.Item(4).ColumnWidth = $src.Columns.Item(4).ColumnWidth
$dst.Columns.Item(6).ColumnWidth = $src.Columns.Item(6).ColumnWidth
$dst.Columns.Item(8).ColumnWidth = $src.Columns.Item(8).ColumnWidth

$dst.Rows.Item(1).RowHeight = $src.Rows.Item(1).RowHeight
$dst.Rows.Item(2).RowHeight = $src.Rows.Item(2).RowHeight
$dst.Rows.Item(3).RowHeight = $src.Rows.Item(3).RowHeight
$dst.Rows.Item(4).RowHeight = $src.Rows.Item(4).RowHeight
$dst.Rows.Item(5).RowHeight = 126
$dst.Rows.Item(6).RowHeight = 108
$dst.Rows.Item(7).RowHeight = $src.Rows.Item(7).RowHeight

# ------------------------------------------------------------------
# 5) Selection / active sheet bookkeeping: the author ended on the
#    "Text Case" tab with F7:G7 selected, having left "Replace
#    Substrings" on a whole-column selection.
# ------------------------------------------------------------------
$src.Range("A1").EntireColumn.Select()
$dst.Activate()
$dst.Range("F7:G7").Select()
